# Append 9 new workout-log rows (rows 281-289) to the GymWorkouts sheet,
# reusing the same column layout / number formats as the existing data,
# and introduces one new Exercise Name value ("Seated Row").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    # A(Id) B(DateId) C(Date-serial) D(Week) E(Month) F(Year) G(Day)      H(Exercise)                          I(Weight) J(Sets) K(Reps)
    @(280, 34, 43093, 51, "December", 2017, "Sunday", "Bicep Curl",                         50, 4,  8),
    @(281, 34, 43093, 51, "December", 2017, "Sunday", "Tricep Pull down",                   40, 4,  8),
    @(282, 34, 43093, 51, "December", 2017, "Sunday", "Shoulder Press",                      24, 4,  8),
    @(283, 34, 43093, 51, "December", 2017, "Sunday", "Shoulder Shrug",                      24, 4,  8),
    @(284, 34, 43093, 51, "December", 2017, "Sunday", "Lying Bicep Curl",                    34, 4,  8),
    @(285, 34, 43093, 51, "December", 2017, "Sunday", "Seated Row",                          61, 4,  8),
    @(286, 34, 43093, 51, "December", 2017, "Sunday", "V-up crunches with medicine ball",     8, 3, 10),
    @(287, 34, 43093, 51, "December", 2017, "Sunday", "Barbell twists",                       0, 3, 10),
    @(288, 34, 43093, 51, "December", 2017, "Sunday", "Russian Twists",                      10, 3, 20)
)

$startRow = 281
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
    $ws.Range("I$r").Value = $row[8]
    $ws.Range("J$r").Value = $row[9]
    $ws.Range("K$r").Value = $row[10]
}

# Reflect the new selection from the source workbook (scrolled/selected
# cell after the edit). The freeze-pane split itself (top row frozen)
# is unchanged.
$ws.Range("H292").Select()
